$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20241004-083336-"

# The reference date (column G, "Dt. Referencia") moved forward one day
# for every data row (45568 -> 45569, i.e. 2024-10-03 -> 2024-10-04)
$ws.Range("G2:G274").Value = 45569

# Row 108: Vl. Projetado (D) and Vl. Total (H) were updated
$ws.Range("D108").Value = 70037.11
$ws.Range("H108").Value = 70476.36

# Row 161: Saldo Previsto (E) and Vl. Total (H) were updated
$ws.Range("E161").Value = 292.08999999999997
$ws.Range("H161").Value = 292.08999999999997
